$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = '$EXECUTION_REQUEST_NUMBER()example@example.com'
$ws.Range("C4:D4").Select()
